$wb = $excel.ActiveWorkbook
$wsRush = $wb.Worksheets.Item("Rushing")
$wsRecv = $wb.Worksheets.Item("Receiving")

# ---------------------------------------------------------------------------
# Rushing sheet (Week 15 logged stats corrections)
# ---------------------------------------------------------------------------

# D.Henderson (row 4)
$wsRush.Range("C4").Value = 95
$wsRush.Range("D4").Value = 49
$wsRush.Range("E4").Value = 20
$wsRush.Range("F4").Value = 27

# S.Michel (row 5)
$wsRush.Range("C5").Value = 60
$wsRush.Range("E5").Value = 14
$wsRush.Range("F5").Value = 23

# M.Sargent (row 8)
$wsRush.Range("C8").Value = 0
$wsRush.Range("D8").Value = 1
$wsRush.Range("E8").Value = 1
$wsRush.Range("F8").Value = 0

# C.Kupp (row 9)
$wsRush.Range("D9").Value = 1

# ---------------------------------------------------------------------------
# Receiving sheet (Week 15 logged stats corrections)
# ---------------------------------------------------------------------------

# D.Henderson (row 2)
$wsRecv.Range("C2").Value = 38
$wsRecv.Range("D2").Value = 27
$wsRecv.Range("E2").Value = 2
$wsRecv.Range("F2").Value = 2
$wsRecv.Range("G2").Value = 9
$wsRecv.Range("H2").Value = 3

# S.Michel (row 3)
$wsRecv.Range("D3").Value = 15
$wsRecv.Range("E3").Value = 1
$wsRecv.Range("F3").Value = 1

# C.Kupp (row 4)
$wsRecv.Range("C4").Value = 112
$wsRecv.Range("D4").Value = 87
$wsRecv.Range("E4").Value = 40
$wsRecv.Range("F4").Value = 22
$wsRecv.Range("G4").Value = 26
$wsRecv.Range("H4").Value = 18

# V.Jefferson (row 5)
$wsRecv.Range("D5").Value = 32
$wsRecv.Range("E5").Value = 22
$wsRecv.Range("F5").Value = 9
$wsRecv.Range("G5").Value = 13

# B.Skowronek (row 6)
$wsRecv.Range("C6").Value = 22
$wsRecv.Range("D6").Value = 7

# O.Beckham (row 7)
$wsRecv.Range("C7").Value = 51
$wsRecv.Range("D7").Value = 45
$wsRecv.Range("E7").Value = 25
$wsRecv.Range("F7").Value = 10
$wsRecv.Range("G7").Value = 19
$wsRecv.Range("H7").Value = 10

# K.Blanton (row 8)
$wsRecv.Range("C8").Value = 3
$wsRecv.Range("D8").Value = 1
$wsRecv.Range("E8").Value = 0
$wsRecv.Range("F8").Value = 0

# ---------------------------------------------------------------------------
# Week 16 simulated: add new receiver B.Hopkins as row 11
# ---------------------------------------------------------------------------

$wsRecv.Range("A11").Value = 9
$wsRecv.Range("B11").Value = "B.Hopkins"
$wsRecv.Range("C11").Value = 1
$wsRecv.Range("D11").Value = 1
$wsRecv.Range("E11").Value = 0
$wsRecv.Range("F11").Value = 0
$wsRecv.Range("G11").Value = 0
$wsRecv.Range("H11").Value = 0

# Match the row-number column's bold/bordered/centered style (same as A10)
$wsRecv.Range("A10").Copy()
$wsRecv.Range("A11").PasteSpecial(-4122)
$wsRecv.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Make "Rushing" the active tab (matches the committed workbook view state)
# ---------------------------------------------------------------------------
$wsRush.Activate()
